$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Checking")
$ws.Range("B2").Value = 57327.05
$ws.Range("C2").Value = 69089.91
$ws.Range("D2").Value = 46505.97
$ws.Range("E2").Value = 68042.99000000001
$ws.Range("F2").Value = 97647.47
$ws.Range("G2").Value = 338613.39
$ws.Range("B3").Value = -11891.18
$ws.Range("C3").Value = -20277.56
$ws.Range("D3").Value = -16132.08
$ws.Range("E3").Value = -18620.19
$ws.Range("F3").Value = -13311.51
$ws.Range("G3").Value = -80232.52
$ws.Range("B4").Value = -28564.63
$ws.Range("C4").Value = -29994.23
$ws.Range("D4").Value = -27644.43
$ws.Range("E4").Value = -29930.08
$ws.Range("F4").Value = -24965.85
$ws.Range("G4").Value = -141099.22
$ws.Range("B5").Value = -196582.08
$ws.Range("C5").Value = -188143.44
$ws.Range("D5").Value = -185434.74
$ws.Range("E5").Value = -205682.74
$ws.Range("F5").Value = -189920.85
$ws.Range("G5").Value = -965763.85
$ws.Range("B6").Value = -18135.67
$ws.Range("C6").Value = -16482.58
$ws.Range("D6").Value = -15629.56
$ws.Range("E6").Value = -19907.67
$ws.Range("F6").Value = -12846.96
$ws.Range("G6").Value = -83002.44
$ws.Range("B7").Value = 81716.10000000001
$ws.Range("C7").Value = 76255.21000000001
$ws.Range("D7").Value = 76187.89999999999
$ws.Range("E7").Value = 68580.99000000001
$ws.Range("F7").Value = 110584.34
$ws.Range("G7").Value = 413324.54
$ws.Range("B8").Value = -116130.41
$ws.Range("C8").Value = -109552.69
$ws.Range("D8").Value = -122146.94
$ws.Range("E8").Value = -137516.7
$ws.Range("F8").Value = -32813.36
$ws.Range("G8").Value = -518160.1

$ws = $wb.Worksheets.Item("Mastercard")
$ws.Range("B2").Value = -21053.09
$ws.Range("C2").Value = -21186.22
$ws.Range("D2").Value = -18976.53
$ws.Range("E2").Value = -16493.13
$ws.Range("F2").Value = -13038.53
$ws.Range("G2").Value = -90747.5
$ws.Range("B3").Value = -15939.19
$ws.Range("C3").Value = -16646
$ws.Range("D3").Value = -21687.32
$ws.Range("E3").Value = -21682.89
$ws.Range("F3").Value = -14987.77
$ws.Range("G3").Value = -90943.17
$ws.Range("B4").Value = -9200.950000000001
$ws.Range("C4").Value = -10266.23
$ws.Range("D4").Value = -8645.75
$ws.Range("E4").Value = -5122.05
$ws.Range("F4").Value = -7230.53
$ws.Range("G4").Value = -40465.51
$ws.Range("B5").Value = -21959.62
$ws.Range("C5").Value = -23666.99
$ws.Range("D5").Value = -16347.16
$ws.Range("E5").Value = -19788.41
$ws.Range("F5").Value = -20591.47
$ws.Range("G5").Value = -102353.65
$ws.Range("B6").Value = -28553.76
$ws.Range("C6").Value = -21922.47
$ws.Range("D6").Value = -27965.68
$ws.Range("E6").Value = -37691.46
$ws.Range("F6").Value = -27734.45
$ws.Range("G6").Value = -143867.82
$ws.Range("B7").Value = 96706.61
$ws.Range("C7").Value = 93687.91
$ws.Range("D7").Value = 93622.44
$ws.Range("E7").Value = 100777.94
$ws.Range("F7").Value = 83582.75
$ws.Range("G7").Value = 468377.65

$ws = $wb.Worksheets.Item("Savings")
$ws.Range("B2").Value = 452802.19
$ws.Range("C2").Value = 439948.04
$ws.Range("D2").Value = 521104.54
$ws.Range("E2").Value = 400785.61
$ws.Range("F2").Value = 439761.29
$ws.Range("G2").Value = 2254401.67
$ws.Range("B3").Value = 9349.23
$ws.Range("C3").Value = 9672.219999999999
$ws.Range("D3").Value = 7688.38
$ws.Range("E3").Value = 8537.84
$ws.Range("F3").Value = 19512.25
$ws.Range("G3").Value = 54759.92
$ws.Range("B4").Value = 462151.42
$ws.Range("C4").Value = 449620.26
$ws.Range("D4").Value = 528792.92
$ws.Range("E4").Value = 409323.45
$ws.Range("F4").Value = 459273.54
$ws.Range("G4").Value = 2309161.59

$ws = $wb.Worksheets.Item("Visa")
$ws.Range("B2").Value = -17862.78
$ws.Range("C2").Value = -21699.55
$ws.Range("D2").Value = -19316.53
$ws.Range("E2").Value = -15460.94
$ws.Range("F2").Value = -15935.42
$ws.Range("G2").Value = -90275.22
$ws.Range("B3").Value = -20413.86
$ws.Range("C3").Value = -13811.38
$ws.Range("D3").Value = -17132.09
$ws.Range("E3").Value = -17263.2
$ws.Range("F3").Value = -19442.3
$ws.Range("G3").Value = -88062.83
$ws.Range("B4").Value = -2665.66
$ws.Range("C4").Value = -10078.34
$ws.Range("D4").Value = -6792.99
$ws.Range("E4").Value = -9886.1
$ws.Range("F4").Value = -10830.72
$ws.Range("G4").Value = -40253.81
$ws.Range("B5").Value = -20200.8
$ws.Range("C5").Value = -16641.67
$ws.Range("D5").Value = -14819.89
$ws.Range("E5").Value = -18004.1
$ws.Range("F5").Value = -18488.93
$ws.Range("G5").Value = -88155.39
$ws.Range("B6").Value = -29383.14
$ws.Range("C6").Value = -22552.37
$ws.Range("D6").Value = -26062.42
$ws.Range("E6").Value = -35752.62
$ws.Range("F6").Value = -22128.48
$ws.Range("G6").Value = -135879.03
$ws.Range("B7").Value = 90526.24000000001
$ws.Range("C7").Value = 84783.31
$ws.Range("D7").Value = 84123.92
$ws.Range("E7").Value = 96366.96000000001
$ws.Range("F7").Value = 86825.85000000001
$ws.Range("G7").Value = 442626.28

